$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = -0.6647498008515189
$ws.Range("J5").Value = 0.450967951930201
$ws.Range("K5").Value = 0.1125813004828595
$ws.Range("L5").Value = 2.513282465284228
